$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert a "FOREGROUND" column between THEME and BUTTON TYPE.
# The existing D1 ("BUTTON TYPE") becomes "FOREGROUND" and the old
# "BUTTON TYPE" header moves out to the new E1 cell.
$ws.Range("E1").Value2 = $ws.Range("D1").Value2
$ws.Range("D1").Value2 = "FOREGROUND"

# Clear out the old sample data row (A2:B2 held "test" / "20").
$ws.Range("A2:B2").ClearContents()

# Add a new data row further down the sheet holding the lower-cased
# versions of the theme / foreground / button-type choices.
$ws.Range("C28").Value2 = "theme"
$ws.Range("D28").Value2 = "foreground"
$ws.Range("E28").Value2 = "button type"

# Leave the selection where the user last clicked before saving.
$ws.Range("F24").Select() | Out-Null
